{"js": "// Correct the devis (quote) amounts: adaptation des erreurs de calcul.\n// Each pair is [old exact text, new exact text] for a standalone w:t run.\nconst replacements = [\n  [\"5095,00\u20ac\", \"5620,00\u20ac\"],\n  [\"11395,00\u20ac\", \"11020,00\u20ac\"],\n  [\"7560,00\u20ac\", \"6570,00\u20ac\"],\n  [\"28285,00\", \"27445,00\"],\n  [\"54145,00\", \"53305,00\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Correction des montants du devis (adaptation des erreurs de calcul).\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"5095,00\u20ac\";  New = \"5620,00\u20ac\" },\n    @{ Old = \"11395,00\u20ac\"; New = \"11020,00\u20ac\" },\n    @{ Old = \"7560,00\u20ac\";  New = \"6570,00\u20ac\" },\n    @{ Old = \"28285,00\";  New = \"27445,00\" },\n    @{ Old = \"54145,00\";  New = \"53305,00\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, $wdReplaceAll) | Out-Null\n}\n"}
